$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("B").Delete()
try {
  $sortObj = $ws.Sort
  $sortObj.SortFields.Clear()
  $sortObj.SortFields.Add($ws.Range("A1:A23"))
  $sortObj.SetRange($ws.Range("A1:B23"))
  $sortObj.Header = 1
  $sortObj.Apply()
  Write-Host "applied"
} catch {
  Write-Host ("ERROR: " + $_.Exception.Message)
}
